$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that are no longer present in the target data set ---
# Row 28 = "SC 92" and row 26 = "RM 232" in the original layout.
# Delete bottom-most first so the other row index stays valid.
$ws.Rows.Item(28).EntireRow.Delete()
$ws.Rows.Item(26).EntireRow.Delete()

# --- Single cell corrections (rows 2-25 keep their original row numbers) ---
$ws.Cells.Item(2, 6).Value = ""          # F2: 18.03 -> (missing)
$ws.Cells.Item(5, 6).Value = 17.66       # F5: (missing) -> 17.66
$ws.Cells.Item(6, 5).Value = -5.7        # E6: (missing) -> -5.7
$ws.Cells.Item(6, 6).Value = 16.43       # F6: (missing) -> 16.43
$ws.Cells.Item(8, 5).Value = ""          # E8: -6.6 -> (missing)
$ws.Cells.Item(10, 6).Value = ""         # F10: 16.43 -> (missing)
$ws.Cells.Item(12, 5).Value = -5.3       # E12: (missing) -> -5.3
$ws.Cells.Item(13, 6).Value = ""         # F13: 17.1 -> (missing)
$ws.Cells.Item(14, 5).Value = ""         # E14: -5.4 -> (missing)
$ws.Cells.Item(17, 5).Value = -7.3       # E17: (missing) -> -7.3
$ws.Cells.Item(18, 5).Value = -8.5       # E18: (missing) -> -8.5
$ws.Cells.Item(19, 5).Value = ""         # E19: -6.5 -> (missing)
$ws.Cells.Item(20, 5).Value = ""         # E20: -7.2 -> (missing)
$ws.Cells.Item(23, 5).Value = -7         # E23: (missing) -> -7
$ws.Cells.Item(24, 6).Value = 16.78      # F24: (missing) -> 16.78

# --- Corrections on the rows that shifted up after the two deletions ---
# Row 26 "SC 5" unchanged.
$ws.Cells.Item(27, 3).Value = 10         # C27 (SC 101): (missing) -> 10
$ws.Cells.Item(27, 5).Value = ""         # E27 (SC 101): -10 -> (missing)
$ws.Cells.Item(28, 6).Value = ""         # F28 (SC 105): 17.44 -> (missing)
$ws.Cells.Item(29, 3).Value = ""         # C29 (SC 119): 11.2 -> (missing)
$ws.Cells.Item(30, 6).Value = 16.89      # F30 (SC 120): (missing) -> 16.89
# Row 31 "SC 132" unchanged.
$ws.Cells.Item(32, 3).Value = ""         # C32 (SC 193): 10.5 -> (missing)
# Row 33 "SC 232" unchanged.
